$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The E1/F1 column headers were mislabeled - "Fixed Assets" and
# "Current Liabilities" were swapped. Fix the header labels.
$ws.Range("E1").Value = "Current Liabilities"
$ws.Range("F1").Value = "Fixed Assets"

# Correct the Fixed Assets amounts (now correctly under column F).
$ws.Range("F2").Value = 60000
$ws.Range("F3").Value = 61000
$ws.Range("F4").Value = 60000

# Update the active selection to match.
$ws.Range("F5").Select()
